$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 236: Count Number of Trapezoids II (no Tags / column C) ---
$ws.Range("A234:B234").Copy($ws.Range("A236:B236"))
$ws.Range("D234:I234").Copy($ws.Range("D236:I236"))
$ws.Cells.Item(236,1).Value = 3625
$ws.Cells.Item(236,2).Value = "Count Number of Trapezoids II"
$ws.Cells.Item(236,4).Value = "hard"
$ws.Cells.Item(236,5).Value = 0
$ws.Cells.Item(236,6).Value = 1
$ws.Cells.Item(236,7).Value = "55??"
$ws.Cells.Item(236,8).Value = 45994
$ws.Cells.Item(236,9).Value = 45994
$ws.Rows.Item(236).RowHeight = 34

# --- Row 237: Valid Anagram ---
$ws.Range("A235:C235").Copy($ws.Range("A237:C237"))
$ws.Range("D235:I235").Copy($ws.Range("D237:I237"))
$ws.Cells.Item(237,1).Value = 242
$ws.Cells.Item(237,2).Value = "Valid Anagram"
$ws.Cells.Item(237,3).Value = "#hash-table "
$ws.Cells.Item(237,4).Value = "easy"
$ws.Cells.Item(237,5).Value = 1
$ws.Cells.Item(237,6).Value = 0
$ws.Cells.Item(237,7).Value = 5
$ws.Cells.Item(237,8).Value = 45994
$ws.Cells.Item(237,9).Value = 45994
$ws.Rows.Item(237).RowHeight = 17

# --- Row 238: Count Collisions on a Road ---
$ws.Range("A235:C235").Copy($ws.Range("A238:C238"))
$ws.Range("D235:I235").Copy($ws.Range("D238:I238"))
$ws.Cells.Item(238,1).Value = 2211
$ws.Cells.Item(238,2).Value = "Count Collisions on a Road"
$ws.Cells.Item(238,3).Value = "#string #array #simulation "
$ws.Cells.Item(238,4).Value = "medium"
$ws.Cells.Item(238,5).Value = 0
$ws.Cells.Item(238,6).Value = 1
$ws.Cells.Item(238,7).Value = 45
$ws.Cells.Item(238,8).Value = 45995
$ws.Cells.Item(238,9).Value = 45995
$ws.Rows.Item(238).RowHeight = 34

# --- Row 239: Longest Repeating Character Replacement ---
$ws.Range("A235:C235").Copy($ws.Range("A239:C239"))
$ws.Range("D235:I235").Copy($ws.Range("D239:I239"))
$ws.Cells.Item(239,1).Value = 424
$ws.Cells.Item(239,2).Value = "Longest Repeating Character Replacement"
$ws.Cells.Item(239,3).Value = "#two-pointers #sliding-window "
$ws.Cells.Item(239,4).Value = "medium"
$ws.Cells.Item(239,5).Value = 0
$ws.Cells.Item(239,6).Value = 1
$ws.Cells.Item(239,7).Value = "75??"
$ws.Cells.Item(239,8).Value = 45995
$ws.Cells.Item(239,9).Value = 45995
$ws.Rows.Item(239).RowHeight = 34

$ws.Range("H239:I239").Select()

Write-Output "done"
